$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = 69571
$ws.Range("E9").Value = 191359921

$ws.Range("C17").Value = 134741
$ws.Range("E17").Value = 296792831

$ws.Range("C122").Value = 9694
$ws.Range("E122").Value = 31947576

$ws.Range("C164").Value = 50573
$ws.Range("E164").Value = 168517372

$ws.Range("C168").Value = 284958
$ws.Range("E168").Value = 1209031755

$ws.Range("C169").Value = 562583
$ws.Range("E169").Value = 1284389559

$ws.Range("C170").Value = 367327
$ws.Range("D170").Value = 38110
$ws.Range("E170").Value = 2844804389

$ws.Range("C171").Value = 115128
$ws.Range("D171").Value = 20264
$ws.Range("E171").Value = 445409607

$ws.Range("C173").Value = 54387
$ws.Range("E173").Value = 151858805

$ws.Range("C174").Value = 357191
$ws.Range("E174").Value = 1016884274

$ws.Range("C175").Value = 125524
$ws.Range("E175").Value = 812307519

$ws.Range("C179").Value = 235683
$ws.Range("E179").Value = 812488760

$ws.Range("C204").Value = 4757
$ws.Range("E204").Value = 11756409

$ws.Range("C205").Value = 11125
$ws.Range("E205").Value = 44114295

$ws.Range("C209").Value = 5364
$ws.Range("E209").Value = 12211202

$ws.Range("C247").Value = 29424
$ws.Range("E247").Value = 99451778

$ws.Range("C264").Value = 47471
$ws.Range("E264").Value = 81946515
